$wb = $excel.ActiveWorkbook

# --- Rename "Uncut Sheet" to "Uncut_Sheet" and fix the Print_Area defined name ---
$ws = $wb.Worksheets.Item("Uncut Sheet")
$ws.Name = "Uncut_Sheet"

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Uncut_Sheet!Print_Area") {
        $n.RefersTo = "=Uncut_Sheet!`$A`$1:`$G`$42"
    }
}

# --- Add new data rows (expiry / lot tracking) on the Uncut_Sheet ---

# Copy C14's formatting down into C15:C16 so the border matches the filled rows
$ws.Range("C14").Copy()
$ws.Range("C15:C16").PasteSpecial(-4122)

$ws.Range("A14").Value = "15.08.2022"
$ws.Range("B14").Value = 115
$ws.Range("D14").Value = 20
$ws.Range("F14").Formula = "=B14-D14"
$ws.Range("G14").Value = "2024-07"

$ws.Range("A15").Value = "16.08.2022"
$ws.Range("B15").Formula = "=F14"
$ws.Range("D15").Value = 15
$ws.Range("F15").Formula = "=B15-D15"

$ws.Range("B16").Formula = "=F15"

# --- Update the active sheet / selection to match the new workflow focus ---
$ws.Activate()
$ws.Range("B10:C11").Select()
